$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.521.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.898.84'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.894.40'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.40%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -1.40%  '
$ws.Range("E10").Value = '  -1.72%  '
$ws.Range("E11").Value = '  +1.62%  '
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.541.96'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.890.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.551.88'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.23'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("E20").Value = '  -2.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '487.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000166'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +9.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  -1.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.044.37'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.46%  '
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.74'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.847.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.64%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("E37").Value = '  +1.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.139'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.18'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.57%  '
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("E42").Value = '  -2.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '428.92'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.21%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.66%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '142.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.808.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0353'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.58%  '
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.07'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.69%  '
